# Update "想去人数" (number of interested attendees) values on the
# "展览" and "全部类型" sheets to reflect newly generated data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 263
    $ws.Range("F4").Value = 166
    $ws.Range("F5").Value = 12
}
